$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row of data for Julia
$ws.Range("A5").Value = "Julia"
$ws.Range("B5").Value = "04A"
$ws.Range("C5").Value = 2

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0"

# Apply text number format to specific numeric cells
$ws.Range("C2").NumberFormat = "@"
$ws.Range("C5").NumberFormat = "@"

$ws.Range("B5").Select()

$ws.PageSetup.Orientation = 1
